$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.505.69"
$ws.Range("E2").Value = "'  -0.26%  "
$ws.Range("D3").Value = "'1.870.64"
$ws.Range("E3").Value = "'  -0.56%  "
$ws.Range("E4").Value = "'  -1.69%  "
$ws.Range("D5").Value = "'315.51"
$ws.Range("E5").Value = "'  -1.15%  "
$ws.Range("E6").Value = "'  -1.75%  "
$ws.Range("D7").Value = "'0.5080"
$ws.Range("E7").Value = "'  -1.59%  "
$ws.Range("E8").Value = "'  -1.97%  "
$ws.Range("D9").Value = "'0.08364"
$ws.Range("E9").Value = "'  -0.19%  "
$ws.Range("D10").Value = "'42.28"
$ws.Range("E10").Value = "'  +0.00%  "
$ws.Range("D11").Value = "'1.105"
$ws.Range("E11").Value = "'  -0.99%  "
$ws.Range("D12").Value = "'6.200"
$ws.Range("E12").Value = "'  -1.25%  "
$ws.Range("D13").Value = "'1.871.07"
$ws.Range("E13").Value = "'  +2.76%  "
$ws.Range("D14").Value = "'20.37"
$ws.Range("E14").Value = "'  -1.13%  "
$ws.Range("D15").Value = "'7.250"
$ws.Range("E15").Value = "'  -0.08%  "
$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "'  -1.79%  "
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = "'  -1.07%  "
$ws.Range("D18").Value = "'91.32"
$ws.Range("E18").Value = "'  -0.36%  "
$ws.Range("D19").Value = "'0.06731"
$ws.Range("E19").Value = "'  -0.97%  "
$ws.Range("D20").Value = "'17.66"
$ws.Range("E20").Value = "'  -0.86%  "
$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = "'  -1.87%  "
$ws.Range("D22").Value = "'5.908"
$ws.Range("E22").Value = "'  -1.39%  "
$ws.Range("D23").Value = "'28.560.77"
$ws.Range("E23").Value = "'  -0.38%  "
$ws.Range("E24").Value = "'  -1.02%  "
$ws.Range("D25").Value = "'2.207"
$ws.Range("E25").Value = "'  -3.94%  "
$ws.Range("D26").Value = "'2.079.06"
$ws.Range("E26").Value = "'  +2.19%  "
$ws.Range("D27").Value = "'156.84"
$ws.Range("E27").Value = "'  -3.70%  "
$ws.Range("D28").Value = "'20.57"
$ws.Range("E28").Value = "'  -1.52%  "
$ws.Range("D29").Value = "'2.419"
$ws.Range("E29").Value = "'  +2.04%  "
$ws.Range("D30").Value = "'126.04"
$ws.Range("E30").Value = "'  -1.60%  "
$ws.Range("E31").Value = "'  -1.64%  "
$ws.Range("E32").Value = "'  -0.47%  "
$ws.Range("E33").Value = "'  -1.62%  "
$ws.Range("D34").Value = "'3.623"
$ws.Range("E34").Value = "'  -0.61%  "
$ws.Range("D35").Value = "'0.02452"
$ws.Range("E35").Value = "'  +0.20%  "
$ws.Range("D36").Value = "'0.06602"
$ws.Range("E36").Value = "'  +1.05%  "
$ws.Range("D37").Value = "'8.981"
$ws.Range("E37").Value = "'  +0.40%  "
$ws.Range("D38").Value = "'0.2163"
$ws.Range("E38").Value = "'  -1.50%  "
$ws.Range("D39").Value = "'5.042"
$ws.Range("E39").Value = "'  -0.10%  "
$ws.Range("E40").Value = "'  -0.94%  "
$ws.Range("D41").Value = "'1.235"
$ws.Range("E41").Value = "'  -4.35%  "
$ws.Range("D42").Value = "'0.6365"
$ws.Range("E42").Value = "'  -1.50%  "
$ws.Range("E43").Value = "'  -1.77%  "
$ws.Range("E44").Value = "'  -1.66%  "
$ws.Range("D45").Value = "'0.6004"
$ws.Range("E45").Value = "'  -0.84%  "
$ws.Range("D46").Value = "'12.98"
$ws.Range("E46").Value = "'  -1.42%  "
$ws.Range("E47").Value = "'  -1.67%  "
$ws.Range("D48").Value = "'2.001"
$ws.Range("E48").Value = "'  +0.02%  "
$ws.Range("E49").Value = "'  -0.02%  "
$ws.Range("D50").Value = "'122.42"
$ws.Range("E50").Value = "'  +0.14%  "
$ws.Range("D51").Value = "'1.120"
$ws.Range("E51").Value = "'  -8.83%  "
